# Applies the "Updated cryptos list" GitHub Actions data refresh to the
# cryptos sheet: refreshed Price (column D) and Volume(1h) (column E)
# values, plus a few rows whose coin ordering changed (rows 33/34 and
# 46/47/48), which required rewriting Coin/Link/Price/Volume together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("D2").NumberFormat = "@"
    $ws.Range("D2").Value = "69.349.12"
    $ws.Range("E2").Value = "  -0.10%  "
    # Row 3
    $ws.Range("D3").NumberFormat = "@"
    $ws.Range("D3").Value = "3.689.80"
    $ws.Range("E3").Value = "  +0.12%  "
    # Row 5
    $ws.Range("D5").NumberFormat = "@"
    $ws.Range("D5").Value = "679.20"
    $ws.Range("E5").Value = "  -1.04%  "
    # Row 6
    $ws.Range("D6").NumberFormat = "@"
    $ws.Range("D6").Value = "159.30"
    $ws.Range("E6").Value = "  -1.32%  "
    # Row 7
    $ws.Range("E7").Value = "  -0.02%  "
    # Row 8
    $ws.Range("D8").NumberFormat = "@"
    $ws.Range("D8").Value = "0.494"
    $ws.Range("E8").Value = "  -0.24%  "
    # Row 9
    $ws.Range("E9").Value = "  -0.50%  "
    # Row 10
    $ws.Range("D10").NumberFormat = "@"
    $ws.Range("D10").Value = "7.18"
    $ws.Range("E10").Value = "  -1.98%  "
    # Row 11
    $ws.Range("E11").Value = "  +1.14%  "
    # Row 12
    $ws.Range("E12").Value = "  -1.55%  "
    # Row 13
    $ws.Range("D13").NumberFormat = "@"
    $ws.Range("D13").Value = "4.310.55"
    $ws.Range("E13").Value = "  +0.11%  "
    # Row 14
    $ws.Range("E14").Value = "  -1.62%  "
    # Row 15
    $ws.Range("D15").NumberFormat = "@"
    $ws.Range("D15").Value = "3.685.04"
    $ws.Range("E15").Value = "  -0.03%  "
    # Row 16
    $ws.Range("D16").NumberFormat = "@"
    $ws.Range("D16").Value = "69.368.68"
    $ws.Range("E16").Value = "  -0.12%  "
    # Row 17
    $ws.Range("E17").Value = "  +3.00%  "
    # Row 18
    $ws.Range("E18").Value = "  +0.26%  "
    # Row 19
    $ws.Range("E19").Value = "  -0.23%  "
    # Row 20
    $ws.Range("D20").NumberFormat = "@"
    $ws.Range("D20").Value = "468.83"
    $ws.Range("E20").Value = "  -1.13%  "
    # Row 21
    $ws.Range("E21").Value = "  -0.44%  "
    # Row 22
    $ws.Range("E22").Value = "  -0.32%  "
    # Row 23
    $ws.Range("D23").NumberFormat = "@"
    $ws.Range("D23").Value = "80.01"
    $ws.Range("E23").Value = "  +0.33%  "
    # Row 24
    $ws.Range("D24").NumberFormat = "@"
    $ws.Range("D24").Value = "3.834.73"
    $ws.Range("E24").Value = "  +0.08%  "
    # Row 25
    $ws.Range("E25").Value = "  +0.04%  "
    # Row 26
    $ws.Range("E26").Value = "  -4.02%  "
    # Row 27
    $ws.Range("D27").NumberFormat = "@"
    $ws.Range("D27").Value = "10.89"
    $ws.Range("E27").Value = "  -2.66%  "
    # Row 28
    $ws.Range("D28").NumberFormat = "@"
    $ws.Range("D28").Value = "9.16"
    $ws.Range("E28").Value = "  -0.55%  "
    # Row 29
    $ws.Range("D29").NumberFormat = "@"
    $ws.Range("D29").Value = "2.70"
    $ws.Range("E29").Value = "  +0.18%  "
    # Row 30
    $ws.Range("D30").NumberFormat = "@"
    $ws.Range("D30").Value = "1.73"
    $ws.Range("E30").Value = "  -2.61%  "
    # Row 31
    $ws.Range("E31").Value = "  -2.62%  "
    # Row 32
    $ws.Range("E32").Value = "  +0.35%  "
    # Row 33
    $ws.Range("B33").Value = "EthereumClassic"
    $ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    $ws.Range("D33").NumberFormat = "@"
    $ws.Range("D33").Value = "26.95"
    $ws.Range("E33").Value = "  +0.77%  "
    # Row 34
    $ws.Range("B34").Value = "ImmutableX"
    $ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    $ws.Range("D34").NumberFormat = "@"
    $ws.Range("D34").Value = "1.99"
    $ws.Range("E34").Value = "  -2.44%  "
    # Row 35
    $ws.Range("D35").NumberFormat = "@"
    $ws.Range("D35").Value = "3.678.83"
    $ws.Range("E35").Value = "  +0.71%  "
    # Row 36
    $ws.Range("D36").NumberFormat = "@"
    $ws.Range("D36").Value = "0.159"
    $ws.Range("E36").Value = "  -2.89%  "
    # Row 37
    $ws.Range("D37").NumberFormat = "@"
    $ws.Range("D37").Value = "8.32"
    $ws.Range("E37").Value = "  +0.37%  "
    # Row 38
    $ws.Range("D38").NumberFormat = "@"
    $ws.Range("D38").Value = "6.28"
    $ws.Range("E38").Value = "  +1.35%  "
    # Row 40
    $ws.Range("E40").Value = "  -2.00%  "
    # Row 41
    $ws.Range("E41").Value = "  -0.16%  "
    # Row 42
    $ws.Range("D42").NumberFormat = "@"
    $ws.Range("D42").Value = "0.0905"
    $ws.Range("E42").Value = "  -0.84%  "
    # Row 43
    $ws.Range("D43").NumberFormat = "@"
    $ws.Range("D43").Value = "169.44"
    $ws.Range("E43").Value = "  +3.34%  "
    # Row 44
    $ws.Range("D44").NumberFormat = "@"
    $ws.Range("D44").Value = "0.941"
    $ws.Range("E44").Value = "  -0.83%  "
    # Row 45
    $ws.Range("D45").NumberFormat = "@"
    $ws.Range("D45").Value = "47.16"
    $ws.Range("E45").Value = "  -2.34%  "
    # Row 46
    $ws.Range("B46").Value = "FLOKI"
    $ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
    $ws.Range("D46").NumberFormat = "@"
    $ws.Range("D46").Value = "0.000280"
    $ws.Range("E46").Value = "  +0.67%  "
    # Row 47
    $ws.Range("B47").Value = "InjectiveProtocol"
    $ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    $ws.Range("D47").NumberFormat = "@"
    $ws.Range("D47").Value = "28.11"
    $ws.Range("E47").Value = "  -6.64%  "
    # Row 48
    $ws.Range("B48").Value = "dogwifhat"
    $ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
    $ws.Range("D48").NumberFormat = "@"
    $ws.Range("D48").Value = "2.70"
    $ws.Range("E48").Value = "  -0.75%  "
    # Row 49
    $ws.Range("E49").Value = "  +0.01%  "
    # Row 50
    $ws.Range("E50").Value = "  -1.76%  "
    # Row 51
    $ws.Range("D51").NumberFormat = "@"
    $ws.Range("D51").Value = "7.87"
    $ws.Range("E51").Value = "  -1.29%  "
